$wb = $excel.ActiveWorkbook

# Sheet ALC, row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 99
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 99
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 99
$ws.Range("N12").Value = -439
$ws.Range("M12").ClearContents()

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3097.8572
$ws.Range("I98").Value = 2647.5
$ws.Range("J98").Value = 5800
$ws.Range("K98").Value = 2647.5
$ws.Range("L98").Value = 5800
$ws.Range("M98").Value = -1149.5
$ws.Range("N98").Value = -8796

# Sheet ALC, row 99
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 265.75
$ws.Range("I99").Value = 260.85715
$ws.Range("J99").Value = 300
$ws.Range("K99").Value = 782.5714499999999
$ws.Range("L99").Value = 900
$ws.Range("M99").Value = 715.4285500000001
$ws.Range("N99").Value = -3896

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3097.8572
$ws.Range("I122").Value = 2647.5
$ws.Range("J122").Value = 5800
$ws.Range("K122").Value = 7942.5
$ws.Range("L122").Value = 17400
$ws.Range("M122").Value = -5492.5
$ws.Range("N122").Value = -22300

# Sheet ALC, row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 44999.332
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 44999.332
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 44999.332
$ws.Range("N134").Value = -55139.332

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2337.0605
$ws.Range("I138").Value = 1069.1212
$ws.Range("J138").Value = 2971.0303
$ws.Range("K138").Value = 3207.3636
$ws.Range("L138").Value = 8913.090899999999
$ws.Range("M138").Value = 1932.6364
$ws.Range("N138").Value = -19193.0909

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1935
$ws.Range("I2").Value = 1888.8889
$ws.Range("J2").Value = 2350
$ws.Range("K2").Value = 1888.8889
$ws.Range("L2").Value = 2350
$ws.Range("M2").Value = -1775.8889
$ws.Range("N2").Value = -2576

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1251.4
$ws.Range("I61").Value = 918.3
$ws.Range("J61").Value = 1917.6
$ws.Range("K61").Value = 918.3
$ws.Range("L61").Value = 1917.6
$ws.Range("M61").Value = -706.3
$ws.Range("N61").Value = -2341.6

# Sheet ARM, row 106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 50000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 50000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1935
$ws.Range("I116").Value = 1888.8889
$ws.Range("J116").Value = 2350
$ws.Range("K116").Value = 1888.8889
$ws.Range("L116").Value = 2350
$ws.Range("M116").Value = 405.1111000000001
$ws.Range("N116").Value = -6938

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1148.2667
$ws.Range("I122").Value = 1086.4615
$ws.Range("J122").Value = 1550
$ws.Range("K122").Value = 3259.3845
$ws.Range("L122").Value = 4650
$ws.Range("M122").Value = -809.3844999999997
$ws.Range("N122").Value = -9550

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1251.4
$ws.Range("I136").Value = 918.3
$ws.Range("J136").Value = 1917.6
$ws.Range("K136").Value = 2754.9
$ws.Range("L136").Value = 5752.799999999999
$ws.Range("M136").Value = -204.8999999999996
$ws.Range("N136").Value = -10852.8

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1935
$ws.Range("I3").Value = 1888.8889
$ws.Range("J3").Value = 2350
$ws.Range("K3").Value = 1888.8889
$ws.Range("L3").Value = 2350
$ws.Range("M3").Value = -1774.8889
$ws.Range("N3").Value = -2578

# Sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 155.61539
$ws.Range("I7").Value = 144.6
$ws.Range("J7").Value = 192.33333
$ws.Range("K7").Value = 144.6
$ws.Range("L7").Value = 192.33333
$ws.Range("M7").Value = -31.59999999999999
$ws.Range("N7").Value = -418.33333

# Sheet CRP, row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 10977.333
$ws.Range("I41").Value = 3666.6667
$ws.Range("J41").Value = 14632.667
$ws.Range("K41").Value = 3666.6667
$ws.Range("L41").Value = 14632.667
$ws.Range("M41").Value = -3238.6667
$ws.Range("N41").Value = -15488.667

# Sheet CRP, row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9407
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 9407
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 9407
$ws.Range("N50").Value = -10657

# Sheet CRP, row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9414.75
$ws.Range("I51").Value = 9295
$ws.Range("J51").Value = 9534.5
$ws.Range("K51").Value = 9295
$ws.Range("L51").Value = 9534.5
$ws.Range("M51").Value = -8559
$ws.Range("N51").Value = -11006.5

# Sheet CRP, row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 14229.25
$ws.Range("I59").Value = 8000
$ws.Range("J59").Value = 16305.667
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 16305.667
$ws.Range("M59").Value = -6855
$ws.Range("N59").Value = -18595.667

# Sheet CRP, row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 7273.6665
$ws.Range("I60").Value = 5050
$ws.Range("J60").Value = 8385.5
$ws.Range("K60").Value = 5050
$ws.Range("L60").Value = 8385.5
$ws.Range("M60").Value = -4539
$ws.Range("N60").Value = -9407.5

# Sheet CRP, row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 9414.75
$ws.Range("I61").Value = 9295
$ws.Range("J61").Value = 9534.5
$ws.Range("K61").Value = 9295
$ws.Range("L61").Value = 9534.5
$ws.Range("M61").Value = -8947
$ws.Range("N61").Value = -10230.5

# Sheet CUL, row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3676645.5
$ws.Range("I2").Value = 6536100
$ws.Range("J2").Value = 203.57143
$ws.Range("K2").Value = 39216600
$ws.Range("L2").Value = 1221.42858
$ws.Range("M2").Value = -39216487
$ws.Range("N2").Value = -1447.42858

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 465.16
$ws.Range("I5").Value = 316.61905
$ws.Range("J5").Value = 1245
$ws.Range("K5").Value = 949.85715
$ws.Range("L5").Value = 3735
$ws.Range("M5").Value = -837.85715
$ws.Range("N5").Value = -3959

# Sheet CUL, row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 48
$ws.Range("I13").Value = 48
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 144
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 24
$ws.Range("N13").ClearContents()

# Sheet CUL, row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 198.8
$ws.Range("I14").Value = 198.8
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 596.4000000000001
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -423.4000000000001

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1111.3
$ws.Range("I122").Value = 650
$ws.Range("J122").Value = 1144.25
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 10298.25
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -15198.25

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5274619.5
$ws.Range("I131").Value = 45546628
$ws.Range("J131").Value = 904.3333
$ws.Range("K131").Value = 136639884
$ws.Range("L131").Value = 2712.9999
$ws.Range("M131").Value = -136634844
$ws.Range("N131").Value = -12792.9999

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 465.16
$ws.Range("I135").Value = 316.61905
$ws.Range("J135").Value = 1245
$ws.Range("K135").Value = 2849.57145
$ws.Range("L135").Value = 11205
$ws.Range("M135").Value = -314.5714500000004
$ws.Range("N135").Value = -16275

# Sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 42810.418
$ws.Range("I139").Value = 44454.348
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 133363.044
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -128223.044
$ws.Range("N139").Value = -25280

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 152717.9
$ws.Range("I140").Value = 178808.17
$ws.Range("J140").Value = 4873
$ws.Range("K140").Value = 536424.51
$ws.Range("L140").Value = 14619
$ws.Range("M140").Value = -531244.51
$ws.Range("N140").Value = -24979

# Sheet GSM, row 63
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Sheet GSM, row 66
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1060.5238
$ws.Range("I102").Value = 1016.875
$ws.Range("J102").Value = 1200.2
$ws.Range("K102").Value = 1016.875
$ws.Range("L102").Value = 1200.2
$ws.Range("M102").Value = 605.125
$ws.Range("N102").Value = -4444.2

# Sheet GSM, row 104
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 40000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 40000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

# Sheet GSM, row 131
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 26750
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 26750
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 26750
$ws.Range("N131").Value = -36830

# Sheet WVR, row 64
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 15332.833
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 15332.833
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 15332.833
$ws.Range("N64").Value = -15828.833

# Sheet WVR, row 67
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 15332.833
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 15332.833
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 15332.833
$ws.Range("N67").Value = -17048.833

# Sheet WVR, row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 58326
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 58326
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 58326
$ws.Range("N135").Value = -68466
